$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 was a "blank" titration row (only C/D/E pre-filled). Fill in the
# new sample taken 5/9/2018 (titrated 24-36): Date, CRM reading, and the
# With Junk / end of sample flags - matching the pattern used by row 24.

# A25: date 5/9/2018 (serial 43229). Copy the date format from A24 so the
# new cell reuses the existing date style instead of minting a new one.
$ws.Range("A25").Value = 43229
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null

# B25: new CRM reading for this sample; D25's shared formula recalculates
# automatically from this.
$ws.Range("B25").Value = 2217.6104390258802

# F25/G25: flag this sample the same way row 24 was flagged.
$ws.Range("F25").Value = "With Junk"
$ws.Range("G25").Value = "end of sample"

# Move the selection down to A26, matching where the user left off entering
# data.
$ws.Range("A26").Select() | Out-Null
